$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 871.871
$ws.Range("J17").Value = 653.3913
$ws.Range("L17").Value = 1960.1739
$ws.Range("N17").Value = -2296.1739
# Row 32
$ws.Range("H32").Value = 1668
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
# Row 40
$ws.Range("H40").Value = 1327.75
$ws.Range("I40").Value = 1400
$ws.Range("J40").Value = 1303.6666
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 1303.6666
$ws.Range("M40").Value = -1225
$ws.Range("N40").Value = -1653.6666
# Row 86
$ws.Range("H86").Value = 1857
$ws.Range("I86").Value = 999.5
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 999.5
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = 123.5
$ws.Range("N86").Value = -4446
# Row 89
$ws.Range("H89").Value = 1857
$ws.Range("I89").Value = 999.5
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 4997.5
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = 618.5
$ws.Range("N89").Value = -22232
# Row 129
$ws.Range("H129").Value = 839.4400000000001
$ws.Range("J129").Value = 899.97754
$ws.Range("L129").Value = 2699.93262
$ws.Range("N129").Value = -12699.93262
# Row 137
$ws.Range("H137").Value = 2536.1904
$ws.Range("I137").Value = 1838.6666
$ws.Range("K137").Value = 5515.9998
$ws.Range("M137").Value = -2965.9998
# Row 141
$ws.Range("H141").Value = 42255.074
$ws.Range("I141").Value = 47116.957
$ws.Range("J141").Value = 3360
$ws.Range("K141").Value = 141350.871
$ws.Range("L141").Value = 10080
$ws.Range("M141").Value = -136170.871
$ws.Range("N141").Value = -20440

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1289.4572
$ws.Range("I61").Value = 1367.625
$ws.Range("K61").Value = 1367.625
$ws.Range("M61").Value = -1155.625
# Row 132
$ws.Range("H132").Value = 2785.125
$ws.Range("I132").Value = 1383.4546
$ws.Range("J132").Value = 3971.1538
$ws.Range("K132").Value = 4150.3638
$ws.Range("L132").Value = 11913.4614
$ws.Range("M132").Value = -1620.3638
$ws.Range("N132").Value = -16973.4614
# Row 136
$ws.Range("H136").Value = 1289.4572
$ws.Range("I136").Value = 1367.625
$ws.Range("K136").Value = 4102.875
$ws.Range("M136").Value = -1552.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 72500
$ws.Range("J92").Value = 72500
$ws.Range("L92").Value = 72500
$ws.Range("N92").Value = -77492
# Row 134
$ws.Range("H134").Value = 2984.2727
$ws.Range("I134").Value = 1843.6
$ws.Range("J134").Value = 5428.5713
$ws.Range("K134").Value = 5530.799999999999
$ws.Range("L134").Value = 16285.7139
$ws.Range("M134").Value = -2995.799999999999
$ws.Range("N134").Value = -21355.7139

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 867.2308
$ws.Range("I22").Value = 300.33334
$ws.Range("J22").Value = 1037.3
$ws.Range("K22").Value = 300.33334
$ws.Range("L22").Value = 1037.3
$ws.Range("M22").Value = 49.66665999999998
$ws.Range("N22").Value = -1737.3
# Row 31
$ws.Range("H31").Value = 3130.9473
$ws.Range("I31").Value = 1147.2222
$ws.Range("J31").Value = 4916.3
$ws.Range("K31").Value = 1147.2222
$ws.Range("L31").Value = 4916.3
$ws.Range("M31").Value = -852.2221999999999
$ws.Range("N31").Value = -5506.3
# Row 34
$ws.Range("H34").Value = 3130.9473
$ws.Range("I34").Value = 1147.2222
$ws.Range("J34").Value = 4916.3
$ws.Range("K34").Value = 1147.2222
$ws.Range("L34").Value = 4916.3
$ws.Range("M34").Value = -945.2221999999999
$ws.Range("N34").Value = -5320.3
# Row 58
$ws.Range("H58").Value = 1837.5161
$ws.Range("I58").Value = 1669.5358
$ws.Range("K58").Value = 1669.5358
$ws.Range("M58").Value = -1466.5358
# Row 122
$ws.Range("H122").Value = 2694.4167
$ws.Range("I122").Value = 1175.8
$ws.Range("K122").Value = 3527.4
$ws.Range("M122").Value = -1077.4
# Row 132
$ws.Range("H132").Value = 2893.1365
$ws.Range("I132").Value = 2044.5714
$ws.Range("J132").Value = 4378.125
$ws.Range("K132").Value = 6133.7142
$ws.Range("L132").Value = 13134.375
$ws.Range("M132").Value = -3603.7142
$ws.Range("N132").Value = -18194.375
# Row 136
$ws.Range("H136").Value = 1837.5161
$ws.Range("I136").Value = 1669.5358
$ws.Range("K136").Value = 5008.607400000001
$ws.Range("M136").Value = -2458.607400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 1000000000
$ws.Range("J37").Value = 1000000000
$ws.Range("L37").Value = 3000000000
$ws.Range("N37").Value = -3000000224
# Row 76
$ws.Range("H76").Value = 3857.1428
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3857.1428
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 11571.4284
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -12337.4284
# Row 79
$ws.Range("H79").Value = 3857.1428
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3857.1428
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 11571.4284
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -14223.4284
# Row 80
$ws.Range("H80").Value = 30714
$ws.Range("J80").Value = 30714
$ws.Range("L80").Value = 92142
$ws.Range("N80").Value = -94014
# Row 81
$ws.Range("H81").Value = 1153.6
$ws.Range("I81").Value = 817
$ws.Range("K81").Value = 2451
$ws.Range("M81").Value = -1328
# Row 82
$ws.Range("H82").Value = 2202.6
$ws.Range("I82").Value = 1006.5
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 3019.5
$ws.Range("L82").Value = 9000
$ws.Range("M82").Value = -2613.5
$ws.Range("N82").Value = -9812
# Row 83
$ws.Range("H83").Value = 30714
$ws.Range("J83").Value = 30714
$ws.Range("L83").Value = 276426
$ws.Range("N83").Value = -285786
# Row 84
$ws.Range("H84").Value = 1153.6
$ws.Range("I84").Value = 817
$ws.Range("K84").Value = 7353
$ws.Range("M84").Value = -1737
# Row 85
$ws.Range("H85").Value = 2202.6
$ws.Range("I85").Value = 1006.5
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 3019.5
$ws.Range("L85").Value = 9000
$ws.Range("M85").Value = -1615.5
$ws.Range("N85").Value = -11808
# Row 113
$ws.Range("H113").Value = 681.1818
$ws.Range("I113").Value = 679.30304
$ws.Range("J113").Value = 686.8182
$ws.Range("K113").Value = 2037.90912
$ws.Range("L113").Value = 2060.4546
$ws.Range("M113").Value = 132.09088
$ws.Range("N113").Value = -6400.4546

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("N26").Value = -20560
# Row 49
$ws.Range("H49").Value = 8758.200000000001
$ws.Range("J49").Value = 10216
$ws.Range("L49").Value = 10216
$ws.Range("N49").Value = -10584
# Row 50
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -20996
# Row 70
$ws.Range("H70").Value = 6159.045
$ws.Range("I70").Value = 5739.268
$ws.Range("J70").Value = 8296.091
$ws.Range("K70").Value = 5739.268
$ws.Range("L70").Value = 8296.091
$ws.Range("M70").Value = -5469.268
$ws.Range("N70").Value = -8836.091
# Row 73
$ws.Range("H73").Value = 6159.045
$ws.Range("I73").Value = 5739.268
$ws.Range("J73").Value = 8296.091
$ws.Range("K73").Value = 5739.268
$ws.Range("L73").Value = 8296.091
$ws.Range("M73").Value = -4803.268
$ws.Range("N73").Value = -10168.091
# Row 80
$ws.Range("H80").Value = 50002620
$ws.Range("I80").Value = 62502530
$ws.Range("K80").Value = 62502530
$ws.Range("M80").Value = -62501532
# Row 83
$ws.Range("H83").Value = 50002620
$ws.Range("I83").Value = 62502530
$ws.Range("K83").Value = 312512650
$ws.Range("M83").Value = -312507658
# Row 102
$ws.Range("H102").Value = 2553.923
$ws.Range("I102").Value = 1620.1
$ws.Range("J102").Value = 5666.6665
$ws.Range("K102").Value = 1620.1
$ws.Range("L102").Value = 5666.6665
$ws.Range("M102").Value = 1.900000000000091
$ws.Range("N102").Value = -8910.666499999999
# Row 111
$ws.Range("H111").Value = 18994.5
$ws.Range("J111").Value = 18994.5
$ws.Range("L111").Value = 18994.5
$ws.Range("N111").Value = -25128.5
# Row 136
$ws.Range("H136").Value = 15584.333
$ws.Range("J136").Value = 15584.333
$ws.Range("L136").Value = 46752.999
$ws.Range("N136").Value = -51852.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 104
$ws.Range("H104").Value = 17939.334
$ws.Range("J104").Value = 17939.334
$ws.Range("L104").Value = 17939.334
$ws.Range("N104").Value = -24927.334
# Row 136
$ws.Range("H136").Value = 3080.1353
$ws.Range("I136").Value = 1099.1765
$ws.Range("J136").Value = 4763.95
$ws.Range("K136").Value = 3297.5295
$ws.Range("L136").Value = 14291.85
$ws.Range("M136").Value = -747.5295000000001
$ws.Range("N136").Value = -19391.85

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3991.0908
$ws.Range("I122").Value = 2627.0667
$ws.Range("K122").Value = 7881.2001
$ws.Range("M122").Value = -5431.2001
# Row 132
$ws.Range("H132").Value = 11112877
$ws.Range("I132").Value = 1112.2727
$ws.Range("J132").Value = 41670230
$ws.Range("K132").Value = 3336.8181
$ws.Range("L132").Value = 125010690
$ws.Range("M132").Value = -806.8181
$ws.Range("N132").Value = -125015750
# Row 136
$ws.Range("H136").Value = 9120
$ws.Range("I136").Value = 6553.1665
$ws.Range("J136").Value = 12200.2
$ws.Range("K136").Value = 19659.4995
$ws.Range("L136").Value = 36600.60000000001
$ws.Range("M136").Value = -17109.4995
$ws.Range("N136").Value = -41700.60000000001
